$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$data = New-Object 'object[,]' 100,2
$data[0,0] = 45812
$data[0,1] = 846085
$data[1,0] = 45813
$data[1,1] = 888671
$data[2,0] = 45817
$data[2,1] = 907427
$data[3,0] = 45818
$data[3,1] = 881779
$data[4,0] = 45819
$data[4,1] = 887470
$data[5,0] = 45820
$data[5,1] = 874461
$data[6,0] = 45821
$data[6,1] = 860946
$data[7,0] = 45824
$data[7,1] = 889311
$data[8,0] = 45825
$data[8,1] = 925023
$data[9,0] = 45826
$data[9,1] = 908175
$data[10,0] = 45827
$data[10,1] = 930072
$data[11,0] = 45828
$data[11,1] = 957046
$data[12,0] = 45831
$data[12,1] = 932764
$data[13,0] = 45832
$data[13,1] = 974838
$data[14,0] = 45833
$data[14,1] = 1085968
$data[15,0] = 45834
$data[15,1] = 1042452
$data[16,0] = 45835
$data[16,1] = 1063568
$data[17,0] = 45838
$data[17,1] = 1058590
$data[18,0] = 45839
$data[18,1] = 1076238
$data[19,0] = 45840
$data[19,1] = 1129606
$data[20,0] = 45841
$data[20,1] = 1091964
$data[21,0] = 45842
$data[21,1] = 1082146
$data[22,0] = 45845
$data[22,1] = 1036585
$data[23,0] = 45846
$data[23,1] = 1057594
$data[24,0] = 45847
$data[24,1] = 1033943
$data[25,0] = 45848
$data[25,1] = 1026565
$data[26,0] = 45849
$data[26,1] = 1026413
$data[27,0] = 45852
$data[27,1] = 1077657
$data[28,0] = 45853
$data[28,1] = 1065374
$data[29,0] = 45854
$data[29,1] = 1029482
$data[30,0] = 45855
$data[30,1] = 1182051
$data[31,0] = 45856
$data[31,1] = 1176308
$data[32,0] = 45859
$data[32,1] = 1174963
$data[33,0] = 45860
$data[33,1] = 1104941
$data[34,0] = 45861
$data[34,1] = 1142894
$data[35,0] = 45862
$data[35,1] = 1119176
$data[36,0] = 45863
$data[36,1] = 1119407
$data[37,0] = 45866
$data[37,1] = 1131330
$data[38,0] = 45867
$data[38,1] = 1116598
$data[39,0] = 45868
$data[39,1] = 1126952
$data[40,0] = 45869
$data[40,1] = 1092961
$data[41,0] = 45870
$data[41,1] = 1063611
$data[42,0] = 45873
$data[42,1] = 1069661
$data[43,0] = 45874
$data[43,1] = 1090836
$data[44,0] = 45875
$data[44,1] = 1091175
$data[45,0] = 45876
$data[45,1] = 1096723
$data[46,0] = 45877
$data[46,1] = 1089483
$data[47,0] = 45880
$data[47,1] = 1089326
$data[48,0] = 45881
$data[48,1] = 1095913
$data[49,0] = 45882
$data[49,1] = 1094951
$data[50,0] = 45883
$data[50,1] = 1103936
$data[51,0] = 45887
$data[51,1] = 1101211
$data[52,0] = 45888
$data[52,1] = 1102799
$data[53,0] = 45889
$data[53,1] = 1110009
$data[54,0] = 45890
$data[54,1] = 1126641
$data[55,0] = 45891
$data[55,1] = 1246325
$data[56,0] = 45894
$data[56,1] = 1255994
$data[57,0] = 45895
$data[57,1] = 1200273
$data[58,0] = 45896
$data[58,1] = 1191348
$data[59,0] = 45897
$data[59,1] = 1224879
$data[60,0] = 45898
$data[60,1] = 1216490
$data[61,0] = 45901
$data[61,1] = 1223392
$data[62,0] = 45902
$data[62,1] = 1265253
$data[63,0] = 45903
$data[63,1] = 1261267
$data[64,0] = 45904
$data[64,1] = 1258094
$data[65,0] = 45905
$data[65,1] = 1262335
$data[66,0] = 45908
$data[66,1] = 1250214
$data[67,0] = 45909
$data[67,1] = 1226883
$data[68,0] = 45910
$data[68,1] = 1181095
$data[69,0] = 45911
$data[69,1] = 1183767
$data[70,0] = 45912
$data[70,1] = 1163265
$data[71,0] = 45915
$data[71,1] = 1155385
$data[72,0] = 45916
$data[72,1] = 1138538
$data[73,0] = 45917
$data[73,1] = 1144070
$data[74,0] = 45918
$data[74,1] = 1124173
$data[75,0] = 45919
$data[75,1] = 1101361
$data[76,0] = 45922
$data[76,1] = 1145804
$data[77,0] = 45923
$data[77,1] = 1140189
$data[78,0] = 45924
$data[78,1] = 1132218
$data[79,0] = 45925
$data[79,1] = 1143621
$data[80,0] = 45926
$data[80,1] = 1094147
$data[81,0] = 45929
$data[81,1] = 1104955
$data[82,0] = 45930
$data[82,1] = 1114247
$data[83,0] = 45931
$data[83,1] = 1150678
$data[84,0] = 45932
$data[84,1] = 1136943
$data[85,0] = 45940
$data[85,1] = 1110940
$data[86,0] = 45943
$data[86,1] = 1087994
$data[87,0] = 45944
$data[87,1] = 1119313
$data[88,0] = 45945
$data[88,1] = 1139829
$data[89,0] = 45946
$data[89,1] = 1238817
$data[90,0] = 45947
$data[90,1] = 1201451
$data[91,0] = 45950
$data[91,1] = 1292434
$data[92,0] = 45951
$data[92,1] = 1320611
$data[93,0] = 45952
$data[93,1] = 1320324
$data[94,0] = 45953
$data[94,1] = 1272469
$data[95,0] = 45954
$data[95,1] = 1275806
$data[96,0] = 45957
$data[96,1] = 1232334
$data[97,0] = 45958
$data[97,1] = 1192738
$data[98,0] = 45959
$data[98,1] = 0
$data[99,0] = 45960
$data[99,1] = 0
$ws.Range("A2:B101").Value = $data

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$data = New-Object 'object[,]' 100,2
$data[0,0] = 45812
$data[0,1] = 900315
$data[1,0] = 45813
$data[1,1] = 820720
$data[2,0] = 45817
$data[2,1] = 816345
$data[3,0] = 45818
$data[3,1] = 835990
$data[4,0] = 45819
$data[4,1] = 839855
$data[5,0] = 45820
$data[5,1] = 850249
$data[6,0] = 45821
$data[6,1] = 844395
$data[7,0] = 45824
$data[7,1] = 854532
$data[8,0] = 45825
$data[8,1] = 819415
$data[9,0] = 45826
$data[9,1] = 843518
$data[10,0] = 45827
$data[10,1] = 842601
$data[11,0] = 45828
$data[11,1] = 804433
$data[12,0] = 45831
$data[12,1] = 799185
$data[13,0] = 45832
$data[13,1] = 821966
$data[14,0] = 45833
$data[14,1] = 794909
$data[15,0] = 45834
$data[15,1] = 791816
$data[16,0] = 45835
$data[16,1] = 742435
$data[17,0] = 45838
$data[17,1] = 785251
$data[18,0] = 45839
$data[18,1] = 767666
$data[19,0] = 45840
$data[19,1] = 767609
$data[20,0] = 45841
$data[20,1] = 758788
$data[21,0] = 45842
$data[21,1] = 745167
$data[22,0] = 45845
$data[22,1] = 752141
$data[23,0] = 45846
$data[23,1] = 803995
$data[24,0] = 45847
$data[24,1] = 824148
$data[25,0] = 45848
$data[25,1] = 827979
$data[26,0] = 45849
$data[26,1] = 830467
$data[27,0] = 45852
$data[27,1] = 816626
$data[28,0] = 45853
$data[28,1] = 822992
$data[29,0] = 45854
$data[29,1] = 813874
$data[30,0] = 45855
$data[30,1] = 902508
$data[31,0] = 45856
$data[31,1] = 873811
$data[32,0] = 45859
$data[32,1] = 876387
$data[33,0] = 45860
$data[33,1] = 893428
$data[34,0] = 45861
$data[34,1] = 822456
$data[35,0] = 45862
$data[35,1] = 838501
$data[36,0] = 45863
$data[36,1] = 823574
$data[37,0] = 45866
$data[37,1] = 826938
$data[38,0] = 45867
$data[38,1] = 733127
$data[39,0] = 45868
$data[39,1] = 733079
$data[40,0] = 45869
$data[40,1] = 749212
$data[41,0] = 45870
$data[41,1] = 704975
$data[42,0] = 45873
$data[42,1] = 736788
$data[43,0] = 45874
$data[43,1] = 767552
$data[44,0] = 45875
$data[44,1] = 762009
$data[45,0] = 45876
$data[45,1] = 753374
$data[46,0] = 45877
$data[46,1] = 745186
$data[47,0] = 45880
$data[47,1] = 755508
$data[48,0] = 45881
$data[48,1] = 742505
$data[49,0] = 45882
$data[49,1] = 775500
$data[50,0] = 45883
$data[50,1] = 764677
$data[51,0] = 45887
$data[51,1] = 781519
$data[52,0] = 45888
$data[52,1] = 796771
$data[53,0] = 45889
$data[53,1] = 773676
$data[54,0] = 45890
$data[54,1] = 793106
$data[55,0] = 45891
$data[55,1] = 781391
$data[56,0] = 45894
$data[56,1] = 788230
$data[57,0] = 45895
$data[57,1] = 797183
$data[58,0] = 45896
$data[58,1] = 865178
$data[59,0] = 45897
$data[59,1] = 825984
$data[60,0] = 45898
$data[60,1] = 830104
$data[61,0] = 45901
$data[61,1] = 807245
$data[62,0] = 45902
$data[62,1] = 826175
$data[63,0] = 45903
$data[63,1] = 831615
$data[64,0] = 45904
$data[64,1] = 838479
$data[65,0] = 45905
$data[65,1] = 844936
$data[66,0] = 45908
$data[66,1] = 852239
$data[67,0] = 45909
$data[67,1] = 826186
$data[68,0] = 45910
$data[68,1] = 796519
$data[69,0] = 45911
$data[69,1] = 825404
$data[70,0] = 45912
$data[70,1] = 821582
$data[71,0] = 45915
$data[71,1] = 818173
$data[72,0] = 45916
$data[72,1] = 807982
$data[73,0] = 45917
$data[73,1] = 810912
$data[74,0] = 45918
$data[74,1] = 810527
$data[75,0] = 45919
$data[75,1] = 812306
$data[76,0] = 45922
$data[76,1] = 779144
$data[77,0] = 45923
$data[77,1] = 774484
$data[78,0] = 45924
$data[78,1] = 801979
$data[79,0] = 45925
$data[79,1] = 805971
$data[80,0] = 45926
$data[80,1] = 755663
$data[81,0] = 45929
$data[81,1] = 766546
$data[82,0] = 45930
$data[82,1] = 788039
$data[83,0] = 45931
$data[83,1] = 790315
$data[84,0] = 45932
$data[84,1] = 780190
$data[85,0] = 45940
$data[85,1] = 808632
$data[86,0] = 45943
$data[86,1] = 803046
$data[87,0] = 45944
$data[87,1] = 773781
$data[88,0] = 45945
$data[88,1] = 805153
$data[89,0] = 45946
$data[89,1] = 846206
$data[90,0] = 45947
$data[90,1] = 856619
$data[91,0] = 45950
$data[91,1] = 859252
$data[92,0] = 45951
$data[92,1] = 896365
$data[93,0] = 45952
$data[93,1] = 902316
$data[94,0] = 45953
$data[94,1] = 892101
$data[95,0] = 45954
$data[95,1] = 930298
$data[96,0] = 45957
$data[96,1] = 1092795
$data[97,0] = 45958
$data[97,1] = 1108744
$data[98,0] = 45959
$data[98,1] = 0
$data[99,0] = 45960
$data[99,1] = 0
$ws.Range("A2:B101").Value = $data

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$data = New-Object 'object[,]' 100,2
$data[0,0] = 45812
$data[0,1] = 798437
$data[1,0] = 45813
$data[1,1] = 802711
$data[2,0] = 45817
$data[2,1] = 804930
$data[3,0] = 45818
$data[3,1] = 814629
$data[4,0] = 45819
$data[4,1] = 809358
$data[5,0] = 45820
$data[5,1] = 829636
$data[6,0] = 45821
$data[6,1] = 848094
$data[7,0] = 45824
$data[7,1] = 861799
$data[8,0] = 45825
$data[8,1] = 905816
$data[9,0] = 45826
$data[9,1] = 883503
$data[10,0] = 45827
$data[10,1] = 901407
$data[11,0] = 45828
$data[11,1] = 931425
$data[12,0] = 45831
$data[12,1] = 1026036
$data[13,0] = 45832
$data[13,1] = 1110078
$data[14,0] = 45833
$data[14,1] = 1052143
$data[15,0] = 45834
$data[15,1] = 989804
$data[16,0] = 45835
$data[16,1] = 954134
$data[17,0] = 45838
$data[17,1] = 934035
$data[18,0] = 45839
$data[18,1] = 873250
$data[19,0] = 45840
$data[19,1] = 842371
$data[20,0] = 45841
$data[20,1] = 808869
$data[21,0] = 45842
$data[21,1] = 774188
$data[22,0] = 45845
$data[22,1] = 781229
$data[23,0] = 45846
$data[23,1] = 809027
$data[24,0] = 45847
$data[24,1] = 829883
$data[25,0] = 45848
$data[25,1] = 837460
$data[26,0] = 45849
$data[26,1] = 815850
$data[27,0] = 45852
$data[27,1] = 810222
$data[28,0] = 45853
$data[28,1] = 781567
$data[29,0] = 45854
$data[29,1] = 804540
$data[30,0] = 45855
$data[30,1] = 800434
$data[31,0] = 45856
$data[31,1] = 794553
$data[32,0] = 45859
$data[32,1] = 831775
$data[33,0] = 45860
$data[33,1] = 830788
$data[34,0] = 45861
$data[34,1] = 808399
$data[35,0] = 45862
$data[35,1] = 854562
$data[36,0] = 45863
$data[36,1] = 856654
$data[37,0] = 45866
$data[37,1] = 900574
$data[38,0] = 45867
$data[38,1] = 976892
$data[39,0] = 45868
$data[39,1] = 941831
$data[40,0] = 45869
$data[40,1] = 1006922
$data[41,0] = 45870
$data[41,1] = 976549
$data[42,0] = 45873
$data[42,1] = 954653
$data[43,0] = 45874
$data[43,1] = 936919
$data[44,0] = 45875
$data[44,1] = 928249
$data[45,0] = 45876
$data[45,1] = 958188
$data[46,0] = 45877
$data[46,1] = 945869
$data[47,0] = 45880
$data[47,1] = 941007
$data[48,0] = 45881
$data[48,1] = 1006990
$data[49,0] = 45882
$data[49,1] = 1037043
$data[50,0] = 45883
$data[50,1] = 1049338
$data[51,0] = 45887
$data[51,1] = 1049168
$data[52,0] = 45888
$data[52,1] = 990011
$data[53,0] = 45889
$data[53,1] = 956613
$data[54,0] = 45890
$data[54,1] = 1011332
$data[55,0] = 45891
$data[55,1] = 1050251
$data[56,0] = 45894
$data[56,1] = 1069181
$data[57,0] = 45895
$data[57,1] = 1026074
$data[58,0] = 45896
$data[58,1] = 1152419
$data[59,0] = 45897
$data[59,1] = 1221324
$data[60,0] = 45898
$data[60,1] = 1414278
$data[61,0] = 45901
$data[61,1] = 1426444
$data[62,0] = 45902
$data[62,1] = 1473132
$data[63,0] = 45903
$data[63,1] = 1444522
$data[64,0] = 45904
$data[64,1] = 1468285
$data[65,0] = 45905
$data[65,1] = 1479565
$data[66,0] = 45908
$data[66,1] = 1485047
$data[67,0] = 45909
$data[67,1] = 1507587
$data[68,0] = 45910
$data[68,1] = 1563474
$data[69,0] = 45911
$data[69,1] = 1592538
$data[70,0] = 45912
$data[70,1] = 1604651
$data[71,0] = 45915
$data[71,1] = 1573322
$data[72,0] = 45916
$data[72,1] = 1607298
$data[73,0] = 45917
$data[73,1] = 1666383
$data[74,0] = 45918
$data[74,1] = 1640537
$data[75,0] = 45919
$data[75,1] = 1648761
$data[76,0] = 45922
$data[76,1] = 1644342
$data[77,0] = 45923
$data[77,1] = 1622619
$data[78,0] = 45924
$data[78,1] = 1637471
$data[79,0] = 45925
$data[79,1] = 1660834
$data[80,0] = 45926
$data[80,1] = 1561614
$data[81,0] = 45929
$data[81,1] = 1562144
$data[82,0] = 45930
$data[82,1] = 1587654
$data[83,0] = 45931
$data[83,1] = 1625506
$data[84,0] = 45932
$data[84,1] = 1565085
$data[85,0] = 45940
$data[85,1] = 1532818
$data[86,0] = 45943
$data[86,1] = 1540761
$data[87,0] = 45944
$data[87,1] = 1516418
$data[88,0] = 45945
$data[88,1] = 1572828
$data[89,0] = 45946
$data[89,1] = 1562795
$data[90,0] = 45947
$data[90,1] = 1685315
$data[91,0] = 45950
$data[91,1] = 1658511
$data[92,0] = 45951
$data[92,1] = 1779600
$data[93,0] = 45952
$data[93,1] = 1806586
$data[94,0] = 45953
$data[94,1] = 1811845
$data[95,0] = 45954
$data[95,1] = 1887798
$data[96,0] = 45957
$data[96,1] = 1965804
$data[97,0] = 45958
$data[97,1] = 1902366
$data[98,0] = 45959
$data[98,1] = 0
$data[99,0] = 45960
$data[99,1] = 0
$ws.Range("A2:B101").Value = $data

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$data = New-Object 'object[,]' 100,2
$data[0,0] = 45812
$data[0,1] = 603423
$data[1,0] = 45813
$data[1,1] = 621714
$data[2,0] = 45817
$data[2,1] = 606346
$data[3,0] = 45818
$data[3,1] = 636158
$data[4,0] = 45819
$data[4,1] = 631383
$data[5,0] = 45820
$data[5,1] = 684619
$data[6,0] = 45821
$data[6,1] = 686183
$data[7,0] = 45824
$data[7,1] = 694219
$data[8,0] = 45825
$data[8,1] = 726255
$data[9,0] = 45826
$data[9,1] = 740388
$data[10,0] = 45827
$data[10,1] = 601383
$data[11,0] = 45828
$data[11,1] = 608959
$data[12,0] = 45831
$data[12,1] = 603980
$data[13,0] = 45832
$data[13,1] = 634012
$data[14,0] = 45833
$data[14,1] = 620416
$data[15,0] = 45834
$data[15,1] = 611138
$data[16,0] = 45835
$data[16,1] = 615295
$data[17,0] = 45838
$data[17,1] = 606541
$data[18,0] = 45839
$data[18,1] = 610583
$data[19,0] = 45840
$data[19,1] = 612116
$data[20,0] = 45841
$data[20,1] = 619249
$data[21,0] = 45842
$data[21,1] = 579939
$data[22,0] = 45845
$data[22,1] = 589345
$data[23,0] = 45846
$data[23,1] = 615315
$data[24,0] = 45847
$data[24,1] = 633102
$data[25,0] = 45848
$data[25,1] = 657498
$data[26,0] = 45849
$data[26,1] = 681323
$data[27,0] = 45852
$data[27,1] = 664185
$data[28,0] = 45853
$data[28,1] = 667876
$data[29,0] = 45854
$data[29,1] = 698537
$data[30,0] = 45855
$data[30,1] = 707470
$data[31,0] = 45856
$data[31,1] = 698981
$data[32,0] = 45859
$data[32,1] = 737261
$data[33,0] = 45860
$data[33,1] = 749481
$data[34,0] = 45861
$data[34,1] = 737252
$data[35,0] = 45862
$data[35,1] = 778003
$data[36,0] = 45863
$data[36,1] = 782972
$data[37,0] = 45866
$data[37,1] = 821664
$data[38,0] = 45867
$data[38,1] = 740027
$data[39,0] = 45868
$data[39,1] = 710025
$data[40,0] = 45869
$data[40,1] = 828931
$data[41,0] = 45870
$data[41,1] = 856342
$data[42,0] = 45873
$data[42,1] = 853303
$data[43,0] = 45874
$data[43,1] = 852631
$data[44,0] = 45875
$data[44,1] = 824998
$data[45,0] = 45876
$data[45,1] = 813968
$data[46,0] = 45877
$data[46,1] = 807931
$data[47,0] = 45880
$data[47,1] = 747899
$data[48,0] = 45881
$data[48,1] = 718148
$data[49,0] = 45882
$data[49,1] = 734096
$data[50,0] = 45883
$data[50,1] = 712303
$data[51,0] = 45887
$data[51,1] = 736661
$data[52,0] = 45888
$data[52,1] = 706693
$data[53,0] = 45889
$data[53,1] = 703749
$data[54,0] = 45890
$data[54,1] = 725211
$data[55,0] = 45891
$data[55,1] = 759451
$data[56,0] = 45894
$data[56,1] = 780558
$data[57,0] = 45895
$data[57,1] = 744264
$data[58,0] = 45896
$data[58,1] = 796040
$data[59,0] = 45897
$data[59,1] = 824994
$data[60,0] = 45898
$data[60,1] = 832119
$data[61,0] = 45901
$data[61,1] = 854001
$data[62,0] = 45902
$data[62,1] = 906667
$data[63,0] = 45903
$data[63,1] = 815396
$data[64,0] = 45904
$data[64,1] = 790047
$data[65,0] = 45905
$data[65,1] = 782861
$data[66,0] = 45908
$data[66,1] = 805492
$data[67,0] = 45909
$data[67,1] = 1302229
$data[68,0] = 45910
$data[68,1] = 1340852
$data[69,0] = 45911
$data[69,1] = 1320774
$data[70,0] = 45912
$data[70,1] = 1246321
$data[71,0] = 45915
$data[71,1] = 1176624
$data[72,0] = 45916
$data[72,1] = 1120876
$data[73,0] = 45917
$data[73,1] = 1124446
$data[74,0] = 45918
$data[74,1] = 1074790
$data[75,0] = 45919
$data[75,1] = 1020586
$data[76,0] = 45922
$data[76,1] = 978509
$data[77,0] = 45923
$data[77,1] = 971206
$data[78,0] = 45924
$data[78,1] = 966777
$data[79,0] = 45925
$data[79,1] = 953411
$data[80,0] = 45926
$data[80,1] = 904018
$data[81,0] = 45929
$data[81,1] = 909924
$data[82,0] = 45930
$data[82,1] = 923281
$data[83,0] = 45931
$data[83,1] = 922552
$data[84,0] = 45932
$data[84,1] = 930346
$data[85,0] = 45940
$data[85,1] = 926251
$data[86,0] = 45943
$data[86,1] = 910097
$data[87,0] = 45944
$data[87,1] = 840682
$data[88,0] = 45945
$data[88,1] = 922184
$data[89,0] = 45946
$data[89,1] = 930526
$data[90,0] = 45947
$data[90,1] = 889583
$data[91,0] = 45950
$data[91,1] = 1001571
$data[92,0] = 45951
$data[92,1] = 975676
$data[93,0] = 45952
$data[93,1] = 995449
$data[94,0] = 45953
$data[94,1] = 1009805
$data[95,0] = 45954
$data[95,1] = 1017010
$data[96,0] = 45957
$data[96,1] = 1053614
$data[97,0] = 45958
$data[97,1] = 961897
$data[98,0] = 45959
$data[98,1] = 0
$data[99,0] = 45960
$data[99,1] = 0
$ws.Range("A2:B101").Value = $data
